# Generate Report for Handoff
#
# The bf804513-be9e-4e4a-bb07-552230a1ed0b.md file has finished translation
# and a fresh handoff xliff was generated, so its status flips from
# "In Translation" to "Ready for handoff" (matching the status already
# shown for fae22e52-bce3-4e87-a08c-5ffa702137a8.md), and the associated
# "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps move forward.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: row 6 is bf804513-be9e-4e4a-bb07-552230a1ed0b.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E6").Value = "Ready for handoff"
$overview.Range("F6").Value = "Ready for handoff"
$overview.Range("G6").Value = "2016-10-20 00:04:19"

# --- "zh-cn" sheet: row 6 is bf804513-be9e-4e4a-bb07-552230a1ed0b.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C6").Value = "Ready for handoff"
$zhcn.Range("H6").Value = "2016-10-20 00:04:08"

# --- "de-de" sheet: row 6 is bf804513-be9e-4e4a-bb07-552230a1ed0b.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C6").Value = "Ready for handoff"
$dede.Range("H6").Value = "2016-10-20 00:04:19"
